# Updated values based on 2.1
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Table -1.1")
$ws2 = $wb.Worksheets.Item("Table - 2.1")
$ws3 = $wb.Worksheets.Item("Table-3.1")

# --- Table -1.1 : update counts / answer text ---
$ws1.Range("C5").Value = 66368
$ws1.Range("C6").Value = 66370
$ws1.Range("C9").Value = "1. Values of key columns have different case`n2. Null values exists in raised_amount_usd column"

# --- Table - 2.1 : update average funding amounts ---
$ws2.Range("C5").Value = 10634054
$ws2.Range("C6").Value = 764564
$ws2.Range("C7").Value = 556607
$ws2.Range("C8").Value = 62111788

# --- Update selections to match final saved state ---
[void]$ws1.Range("C10").Select()
[void]$ws3.Range("C8").Select()

# --- Table-3.1 becomes the active / visible tab, zoomed to 130% ---
[void]$ws3.Activate()
$excel.ActiveWindow.Zoom = 130
